# Colors_Master.xlsx - minor text fixes (casing/wording normalization of
# color names) plus updated selection, per commit "Bunch of minor fixes".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Normalize a handful of color-name labels in column B (case / wording
# fixes - e.g. "Light Nougat" -> "Light nougat").
$ws.Range("B43").Value = "Light nougat"
$ws.Range("B47").Value = "Medium nougat"
$ws.Range("B50").Value = "Dark azur"
$ws.Range("B51").Value = "Medium azur"
$ws.Range("B53").Value = "Medium lavender"
$ws.Range("B55").Value = "Spring yellowish green"
$ws.Range("B57").Value = "Olive green"

# Move the active selection in the frozen (scrolling) pane down to B42.
$ws.Range("B42").Select()
